$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $val
    $ws.Range($cell).Style = "Normal"
}

# --- Row 13 / 14 swap: Solana and WrappedEther swap places ---
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.888.37"
$ws.Range("E13").Value = "  +1.44%  "

$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D14" "20.60"
$ws.Range("E14").Value = "  +0.20%  "

# --- Price (D) and Volume(1h) (E) updates ---
$ws.Range("D2").Value = "28.231.95"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "1.882.16"
$ws.Range("E3").Value = "  +1.29%  "
Set-TextValue "D4" "1.006"
$ws.Range("E4").Value = "  +0.14%  "
Set-TextValue "D5" "315.36"
$ws.Range("E5").Value = "  +1.04%  "
Set-TextValue "D6" "1.008"
$ws.Range("E6").Value = "  +0.49%  "
Set-TextValue "D7" "0.5135"
$ws.Range("E7").Value = "  +1.03%  "
Set-TextValue "D8" "0.3911"
$ws.Range("E8").Value = "  +1.98%  "
Set-TextValue "D9" "0.08382"
$ws.Range("E9").Value = "  +1.70%  "
Set-TextValue "D10" "1.123"
$ws.Range("E10").Value = "  +1.11%  "
Set-TextValue "D11" "41.62"
$ws.Range("E11").Value = "  +0.19%  "
Set-TextValue "D12" "6.232"
$ws.Range("E12").Value = "  +0.45%  "
Set-TextValue "D15" "7.272"
$ws.Range("E15").Value = "  +0.42%  "
Set-TextValue "D16" "1.008"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("E17").Value = "  +0.39%  "
Set-TextValue "D18" "91.17"
$ws.Range("E18").Value = "  +0.47%  "
Set-TextValue "D19" "0.06681"
$ws.Range("E19").Value = "  +0.57%  "
Set-TextValue "D20" "17.80"
$ws.Range("E20").Value = "  +0.82%  "
Set-TextValue "D21" "1.007"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "28.273.72"
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("E24").Value = "  +0.76%  "
Set-TextValue "D25" "2.235"
$ws.Range("E25").Value = "  -0.31%  "
Set-TextValue "D26" "159.81"
$ws.Range("E26").Value = "  +1.57%  "
Set-TextValue "D27" "2.472"
$ws.Range("E27").Value = "  -1.61%  "
Set-TextValue "D28" "20.78"
$ws.Range("E28").Value = "  +1.50%  "
Set-TextValue "D29" "125.41"
$ws.Range("E29").Value = "  +0.51%  "
Set-TextValue "D30" "0.1056"
$ws.Range("E30").Value = "  -0.21%  "
Set-TextValue "D31" "1.043"
$ws.Range("E31").Value = "  +1.03%  "
Set-TextValue "D32" "5.852"
$ws.Range("E32").Value = "  -0.83%  "
Set-TextValue "D33" "3.617"
$ws.Range("E33").Value = "  +0.48%  "
Set-TextValue "D34" "9.576"
$ws.Range("E34").Value = "  +2.20%  "
Set-TextValue "D35" "0.02452"
$ws.Range("E35").Value = "  +1.78%  "
Set-TextValue "D36" "0.06577"
$ws.Range("E36").Value = "  +1.17%  "
Set-TextValue "D37" "0.2217"
$ws.Range("E37").Value = "  +2.04%  "
Set-TextValue "D38" "1.196"
$ws.Range("E38").Value = "  -0.06%  "
Set-TextValue "D39" "0.6487"
$ws.Range("E39").Value = "  -0.73%  "
Set-TextValue "D40" "1.245"
$ws.Range("E40").Value = "  +1.90%  "
Set-TextValue "D41" "4.996"
$ws.Range("E41").Value = "  +0.00%  "
Set-TextValue "D42" "11.22"
$ws.Range("E42").Value = "  +0.47%  "
Set-TextValue "D43" "0.6090"
$ws.Range("E43").Value = "  -0.57%  "
Set-TextValue "D44" "13.02"
$ws.Range("E44").Value = "  -0.90%  "
Set-TextValue "D45" "3.697"
$ws.Range("E45").Value = "  +1.29%  "
Set-TextValue "D46" "1.280"
$ws.Range("E46").Value = "  -0.17%  "
Set-TextValue "D47" "2.015"
$ws.Range("E47").Value = "  +0.17%  "
Set-TextValue "D48" "1.235"
$ws.Range("E48").Value = "  +2.42%  "
Set-TextValue "D49" "121.15"
$ws.Range("E49").Value = "  +0.99%  "
Set-TextValue "D50" "0.06904"
$ws.Range("E50").Value = "  +1.03%  "
Set-TextValue "D51" "77.97"
$ws.Range("E51").Value = "  -0.62%  "
